$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7 and 6 (bottom-up so row indices of remaining rows stay stable)
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# Update changed figures
$ws.Cells.Item(2, 4).Value = 257.173264        # D2

$ws.Cells.Item(3, 2).Value = 44.160058         # B3
$ws.Cells.Item(3, 3).Value = 1                 # C3
$ws.Cells.Item(3, 4).Value = 0.04256           # D3
$ws.Cells.Item(3, 5).Value = 0.83675           # E3

$ws.Cells.Item(4, 2).Value = 225158.462954     # B4
$ws.Cells.Item(4, 3).Value = 217               # C4

$ws.Cells.Item(5, 8).Value = -8.753895999999999 # H5
$ws.Cells.Item(5, 9).Value = 10.800679          # I5
$ws.Cells.Item(5, 10).Value = 0.83675           # J5
